# Update gh-pages to output generated at 456a3b4
#
# The workbook tracks "想去人数" (want-to-go headcount) pulled from
# Bilibili event listings across three sheets:
#   展览 (Exhibitions), 演出 (Performances), 全部类型 (All types / combined)
# A fresh data-generation run bumped the F-column counts for the rows
# below; only column F (numeric headcount) values change, everything
# else (ids, G price, H link, I cover, etc.) stays untouched.

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1222
$ws.Range("F3").Value = 1119
$ws.Range("F4").Value = 872
$ws.Range("F5").Value = 100
$ws.Range("F8").Value = 82
$ws.Range("F9").Value = 41
$ws.Range("F11").Value = 2217
$ws.Range("F12").Value = 1551
$ws.Range("F13").Value = 1235
$ws.Range("F15").Value = 221
$ws.Range("F16").Value = 487
$ws.Range("F17").Value = 713
$ws.Range("F18").Value = 264
$ws.Range("F22").Value = 4176
$ws.Range("F27").Value = 68
$ws.Range("F28").Value = 601
$ws.Range("F30").Value = 60
$ws.Range("F31").Value = 33
$ws.Range("F33").Value = 355
$ws.Range("F34").Value = 918
$ws.Range("F37").Value = 110
$ws.Range("F38").Value = 104

# --- 演出 (Performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 417

# --- 全部类型 (All types / combined) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1222
$ws.Range("F5").Value = 1119
$ws.Range("F6").Value = 872
$ws.Range("F8").Value = 417
$ws.Range("F9").Value = 100
$ws.Range("F12").Value = 82
$ws.Range("F13").Value = 41
$ws.Range("F16").Value = 2217
$ws.Range("F17").Value = 1551
$ws.Range("F18").Value = 1235
$ws.Range("F20").Value = 221
$ws.Range("F21").Value = 487
$ws.Range("F23").Value = 713
$ws.Range("F24").Value = 264
$ws.Range("F28").Value = 4176
$ws.Range("F33").Value = 68
$ws.Range("F34").Value = 601
$ws.Range("F36").Value = 60
$ws.Range("F37").Value = 33
$ws.Range("F39").Value = 355
$ws.Range("F40").Value = 918
$ws.Range("F43").Value = 110
$ws.Range("F44").Value = 104
